$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C19").Value = " - SA 를 파이썬 스크립트를 구동 가능하도록 적용`n - 외부에서 VM으로 메일이 올 수 있도록 port forwarding 공부"

$ws.Range("B12:C12").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
